$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.852.21'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '2.363.62'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("E5").Value = '  -1.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '239.88'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.39'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.601'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.103'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.95'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '36.81'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +13.51%  '
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.27'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.36'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.930'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.97%  '
$ws.Range("D17").Value = '2.370.18'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").Value = '43.817.85'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("E19").Value = '  +1.97%  '
$ws.Range("E20").Value = '  -4.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '77.26'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '253.99'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.88%  '
$ws.Range("E23").Value = '  +3.67%  '
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.88'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -4.71%  '
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.59'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.30'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.38'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.35'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.130'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.38%  '
$ws.Range("E32").Value = '  -1.61%  '
$ws.Range("E33").Value = '  -0.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.46'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.10'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.78%  '
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("E37").Value = '  +4.24%  '
$ws.Range("E38").Value = '  +1.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0280'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.69'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +20.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '20.76'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +8.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '65.05'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +11.71%  '
$ws.Range("E43").Value = '  -4.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.07'
$ws.Range("D44").ClearFormats()
$ws.Range("E45").Value = '  -0.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.53'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.29%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.25'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("B48").Value = 'BinanceUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("E50").Value = '  -2.34%  '
$ws.Range("E51").Value = '  +2.17%  '
